# Hortaliza, Terminal La Palmera de La Serena - Acelga
# "Fruta / hortaliza, semanal" -- weekly append of a new price observation.
#
# The new week's data is inserted as two new rows right above the existing
# row 420 (one "Primera" + one "Segunda" quality record), which pushes every
# row from the old 420..452 down to 422..454 (dimension grows from R452 to
# R454) while keeping all of their values intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows above row 420; everything below (420-452) shifts
# down to 422-454, carrying its existing values/styles with it.
$ws.Rows.Item(420).Resize(2).Insert()

# New row 420: "Primera" quality observation for the new date.
$ws.Cells.Item(420, 1).Value  = 8
$ws.Cells.Item(420, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(420, 3).Value  = "Coquimbo"
$ws.Cells.Item(420, 4).Value  = 44769
$ws.Cells.Item(420, 5).Value  = 4
$ws.Cells.Item(420, 6).Value  = 100112009
$ws.Cells.Item(420, 7).Value  = "Acelga"
$ws.Cells.Item(420, 8).Value  = "Sin especificar"
$ws.Cells.Item(420, 9).Value  = "Primera"
$ws.Cells.Item(420, 10).Value = 2500
$ws.Cells.Item(420, 11).Value = 600
$ws.Cells.Item(420, 12).Value = 700
$ws.Cells.Item(420, 13).Value = 650
$ws.Cells.Item(420, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(420, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(420, 16).Value = 325
$ws.Cells.Item(420, 17).Value = 2
$ws.Cells.Item(420, 18).Value = "Hortaliza"

# New row 421: "Segunda" quality observation for the same new date.
$ws.Cells.Item(421, 1).Value  = 8
$ws.Cells.Item(421, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(421, 3).Value  = "Coquimbo"
$ws.Cells.Item(421, 4).Value  = 44769
$ws.Cells.Item(421, 5).Value  = 4
$ws.Cells.Item(421, 6).Value  = 100112009
$ws.Cells.Item(421, 7).Value  = "Acelga"
$ws.Cells.Item(421, 8).Value  = "Sin especificar"
$ws.Cells.Item(421, 9).Value  = "Segunda"
$ws.Cells.Item(421, 10).Value = 1400
$ws.Cells.Item(421, 11).Value = 500
$ws.Cells.Item(421, 12).Value = 550
$ws.Cells.Item(421, 13).Value = 525
$ws.Cells.Item(421, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(421, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(421, 16).Value = 262
$ws.Cells.Item(421, 17).Value = 2
$ws.Cells.Item(421, 18).Value = "Hortaliza"

# Keep the date column's existing number format (it was carried over by the
# row insert already, but make sure it's explicit/consistent).
$ws.Range("D420:D421").NumberFormat = $ws.Range("D422").NumberFormat
